$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 41, pushing the existing rows 41:51 down to 42:52.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new weekly price record.
$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C41").Value = "Los Lagos"
$ws.Range("D41").Value = 44845
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = 100112012
$ws.Range("G41").Value = "Espinaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 35
$ws.Range("K41").Value = 14000
$ws.Range("L41").Value = 14000
$ws.Range("M41").Value = 14000
$ws.Range("N41").Value = "$/cuna 10 kilos"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 1400
$ws.Range("Q41").Value = 10
$ws.Range("R41").Value = "Hortaliza"
